# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'22.075.78"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -1.17%  "

$ws.Cells.Item(3, 4).Value = "'1.556.77"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.38%  "

$ws.Cells.Item(4, 5).Value = "  +0.25%  "

$ws.Cells.Item(5, 4).Value = "'1.001"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.11%  "

$ws.Cells.Item(6, 4).Value = "'292.14"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.16%  "

$ws.Cells.Item(7, 4).Value = "'0.3973"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +5.12%  "

$ws.Cells.Item(8, 4).Value = "'0.3235"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.50%  "

$ws.Cells.Item(9, 4).Value = "'44.12"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.77%  "

$ws.Cells.Item(10, 4).Value = "'0.07327"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.08%  "

$ws.Cells.Item(11, 4).Value = "'1.085"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -5.76%  "

$ws.Cells.Item(12, 5).Value = "  +0.26%  "

$ws.Cells.Item(13, 4).Value = "'18.96"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -7.31%  "

$ws.Cells.Item(14, 4).Value = "'5.704"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.14%  "

$ws.Cells.Item(15, 4).Value = "'0.00001142"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +5.59%  "

$ws.Cells.Item(16, 4).Value = "'6.654"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.90%  "

$ws.Cells.Item(17, 4).Value = "'1.554.21"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.07%  "

$ws.Cells.Item(18, 4).Value = "'0.06612"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.59%  "

$ws.Cells.Item(19, 4).Value = "'83.77"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -3.01%  "

$ws.Cells.Item(20, 4).Value = "'1.000"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.04%  "

$ws.Cells.Item(21, 4).Value = "'6.320"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.74%  "

$ws.Cells.Item(22, 4).Value = "'15.74"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.04%  "

$ws.Cells.Item(23, 4).Value = "'11.30"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -3.90%  "

$ws.Cells.Item(24, 4).Value = "'22.110.75"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.95%  "

$ws.Cells.Item(25, 4).Value = "'2.352"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.30%  "

$ws.Cells.Item(26, 4).Value = "'2.451"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -6.21%  "

$ws.Cells.Item(27, 4).Value = "'148.50"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.17%  "

$ws.Cells.Item(28, 4).Value = "'18.64"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.98%  "

$ws.Cells.Item(29, 4).Value = "'4.868"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.25%  "

$ws.Cells.Item(30, 4).Value = "'1.733.02"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.44%  "

$ws.Cells.Item(31, 4).Value = "'119.06"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.33%  "

$ws.Cells.Item(32, 5).Value = "  -7.83%  "

$ws.Cells.Item(33, 4).Value = "'5.741"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -3.63%  "

$ws.Cells.Item(34, 4).Value = "'0.08364"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.48%  "

$ws.Cells.Item(35, 4).Value = "'1.625"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -14.93%  "

$ws.Cells.Item(36, 4).Value = "'9.096"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.94%  "

$ws.Cells.Item(37, 4).Value = "'0.02279"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -4.02%  "

$ws.Cells.Item(38, 4).Value = "'0.06142"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -3.05%  "

$ws.Cells.Item(39, 4).Value = "'5.139"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -4.21%  "

$ws.Cells.Item(40, 4).Value = "'1.218"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.88%  "

$ws.Cells.Item(41, 4).Value = "'0.2064"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -4.52%  "

$ws.Cells.Item(42, 4).Value = "'0.9999"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.01%  "

$ws.Cells.Item(43, 4).Value = "'10.78"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -3.11%  "

$ws.Cells.Item(44, 4).Value = "'0.5871"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.96%  "

$ws.Cells.Item(45, 4).Value = "'13.15"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -4.69%  "

$ws.Cells.Item(46, 4).Value = "'3.766"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.29%  "

$ws.Cells.Item(47, 4).Value = "'0.5619"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -5.66%  "

$ws.Cells.Item(48, 4).Value = "'119.07"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.46%  "

$ws.Cells.Item(49, 4).Value = "'1.913"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -4.31%  "

$ws.Cells.Item(50, 4).Value = "'1.141"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.38%  "

$ws.Cells.Item(51, 4).Value = "'0.06851"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -3.59%  "
